# ---------------------------------------------------------------------------
# Applies the "all but glaciers added" commit:
#   * Inserts a new worksheet "AllButGlaciers" right before "Expert"
#     (pushing Expert / transposed ar6 / Frederikse down one slot).
#   * Populates it with the GMSL-rate-minus-glacier-rate comparison table.
#   * Adds the two reference hyperlinks on that sheet.
#   * Updates a couple of stale selections on GMSL / Glaciers sheets.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheet, positioned right before "Expert" -------------

$expertSheet = $wb.Worksheets.Item("Expert")
$ws = $wb.Worksheets.Add($expertSheet)
$ws.Name = "AllButGlaciers"

# --- 2. Header block (rows 1-3), mirrors the PEN / GMSL sheets -------------

$ws.Range("A1").Value = "#"
$ws.Range("A2").Value = "#"
$ws.Range("A3").Value = "#"

$ws.Range("D1").Value = "Frederikse"

$ws.Range("G1").Value = "std2likely"
$ws.Range("H1").Value = 0.95417
$ws.Range("I1").Value = "https://www.earth-syst-sci-data.net/11/1189/2019/s"

$ws.Range("G2").Value = "std290"
$ws.Range("H2").Value = 3.2897
$ws.Range("I2").Value = "https://www.pnas.org/content/114/23/5946"

$ws.Hyperlinks.Add($ws.Range("I1"), "https://www.earth-syst-sci-data.net/11/1189/2019/s") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://www.pnas.org/content/114/23/5946") | Out-Null

$ws.Range("I1").Font.Underline = $true
$ws.Range("I1").Font.Color = 16711680
$ws.Range("I2").Font.Underline = $true

# --- 3. Column headers (row 4) ----------------------------------------------

$ws.Range("A4").Value = "Name"
$ws.Range("B4").Value = "Period start"
$ws.Range("C4").Value = "Period end"
$ws.Range("D4").Value = "FrederikseRate"
$ws.Range("E4").Value = "Rate"
$ws.Range("F4").Value = "RateSigma"
$ws.Range("G4").Value = "Source"
$ws.Range("H4").Value = "Note"
$ws.Range("A4:H4").Font.Bold = $true

$ws.Range("J4").Value = "gmslrate"
$ws.Range("K4").Value = "gmslsigma"
$ws.Range("L4").Value = "glac"
$ws.Range("M4").Value = "glacsigma"

# --- 4. Data rows 5-9: period comparisons -----------------------------------

$names = @("1901-1990", "1971-2018", "1993-2018", "2006-2018", "#1901-2018")
$starts = @(1901, 1971, 1993, 2006, 1901)
$ends = @(1990, 2018, 2018, 2018, 2018)
$gmsl = @(1.35, 2.33, 3.25, 3.69, 1.73)
$gmslSigma = @(0.34653615831230811, 0.47724716539502088, 0.22190473295437274, 0.29181992278931213, 0.27054138675259143)
$glac = @(0.57999999999999996, 0.44, 0.55000000000000004, 0.62, 0.56999999999999995)
$glacSigma = @(0.14590996139465603, 0.13983037966987874, 0.091193725871660011, 0.033437699486275375, 0.13071100708271274)

for ($i = 0; $i -lt 5; $i++) {
    $r = 5 + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $starts[$i]
    $ws.Cells.Item($r, 3).Value = $ends[$i]
    $ws.Cells.Item($r, 5).Formula = "=J$r-L$r"
    $ws.Cells.Item($r, 6).Formula = "=SQRT(K$r^2+ M$r^2)"
    $ws.Cells.Item($r, 7).Value = "AR6 ch9"
    $ws.Cells.Item($r, 10).Value = $gmsl[$i]
    $ws.Cells.Item($r, 11).Value = $gmslSigma[$i]
    $ws.Cells.Item($r, 12).Value = $glac[$i]
    $ws.Cells.Item($r, 13).Value = $glacSigma[$i]
}

# --- 5. Row 10: 1850-1900 baseline ------------------------------------------

$ws.Range("A10").Value = "PI"
$ws.Range("B10").Value = 1850
$ws.Range("C10").Value = 1900
$ws.Range("E10").Formula = "=J10-L10"
$ws.Range("F10").Formula = "=SQRT(K10^2+ M10^2)"
$ws.Range("G10").Value = "Kopp, R. E., Kemp, A. C., Bittermann, K., Horton, B. P., Donnelly, J. P., Gehrels, W. R., Hay, C. C., Mitrovica, J. X., Morrow, E. D., and Rahmstorf, S.: Temperature-driven global sea-level variability in the Common Era, P. Natl. Acad. Sci. USA, 113, E1434-E1441, 2016"
$ws.Range("H10").Value = "1850-1900 Estimate based on this paper provided by Kopp by email "
$ws.Range("J10").Formula = "=0.014*1000/50"
$ws.Range("K10").Formula = "=0.014*1000/50"

# --- 6. Selection state on the new sheet (matches the authored file) -------

$ws.Range("Q11").Select()

# --- 7. Stale-selection touch-ups on other sheets ---------------------------

$gmslSheet = $wb.Worksheets.Item("GMSL")
$gmslSheet.Activate()
$gmslSheet.Range("A1:K10").Select()

$glaciersSheet = $wb.Worksheets.Item("Glaciers")
$glaciersSheet.Activate()
$glaciersSheet.Range("E5:F9").Select()

# --- 8. Re-activate the new sheet (becomes the active tab, like the diff) --

$ws.Activate()
